$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.711.19'
$ws.Range("E2").Value = '  -1.24%  '

$ws.Range("D3").Value = '2.631.01'
$ws.Range("E3").Value = '  +0.23%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.68%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -3.36%  '

$ws.Range("D9").Value = '2.629.49'
$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("E10").Value = '  -3.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.18%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.382'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.02%  '

$ws.Range("E13").Value = '  +0.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("D15").Value = '3.106.41'
$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("E16").Value = '  -2.19%  '

$ws.Range("D17").Value = '63.656.29'
$ws.Range("E17").Value = '  -1.12%  '

$ws.Range("D18").Value = '2.616.72'
$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.29%  '

$ws.Range("E21").Value = '  -3.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '344.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.32%  '

$ws.Range("E23").Value = '  +0.38%  '

$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("E25").Value = '  +9.26%  '

$ws.Range("E26").Value = '  -4.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '605.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.97%  '

$ws.Range("E28").Value = '  -1.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.71%  '

$ws.Range("E30").Value = '  -0.07%  '

$ws.Range("E31").Value = '  -0.68%  '

$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.06'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.75'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.63'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.70%  '

$ws.Range("E37").Value = '  -2.32%  '

$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.74'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.43%  '

$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.06%  '

$ws.Range("E40").Value = '  -2.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '149.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.23%  '

$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("E43").Value = '  +2.54%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.53%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '159.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0587'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.19%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.631'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.71%  '

$ws.Range("E51").Value = '  -0.47%  '

